$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$f = $ws.Range("A1:C1").Font
$f.Name = "Arial"
$f.Size = 10
$f.Family = 3
Write-Output "done"
